$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows for years 2000年-2009年 (original rows 2 through 11),
# which shifts the 2010年-2015年 data (original rows 12-17) up to rows 2-7.
$ws.Range("A2:F11").EntireRow.Delete()
